$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Create new row 83 (mirrors formatting of row 82) and set the date value
#    in column A, matching the pattern of the existing date column.
# ---------------------------------------------------------------------------
$ws.Range("A82").Copy() | Out-Null
$ws.Range("A83").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A83").Value = 46934

# ---------------------------------------------------------------------------
# 2. Add new column BB (header date 45986, plus data down through row 83),
#    mirroring the style of column BA (bold/bordered header in row 1, plain
#    numeric cells below).
# ---------------------------------------------------------------------------
$ws.Range("BA1:BA82").Copy() | Out-Null
$ws.Range("BB1:BB82").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$bbArr = New-Object 'object[,]' 83,1
$bbArr[0,0] = 45986
$bbArr[1,0] = 0.04899821040007168
$bbArr[2,0] = 1.885178963001849
$bbArr[3,0] = -1.307280175329765
$bbArr[4,0] = 2.75267136845396
$bbArr[5,0] = -4.471592960313714
$bbArr[6,0] = -4.942495347902479
$bbArr[7,0] = -5.434114574907241
$bbArr[8,0] = 4.370618442157621
$bbArr[9,0] = -0.9953340506219917
$bbArr[10,0] = 5.692238679293155
$bbArr[11,0] = 5.484876523251003
$bbArr[12,0] = 0.1757137213762547
$bbArr[13,0] = 1.626538719806248
$bbArr[14,0] = 2.437334396728659
$bbArr[15,0] = 1.566144859651857
$bbArr[16,0] = 1.556824096124856
$bbArr[17,0] = 0.21169683169569
$bbArr[18,0] = 0.1582599013804469
$bbArr[19,0] = -0.3253450194449812
$bbArr[20,0] = -0.15853729477206
$bbArr[21,0] = 0.2054976683197225
$bbArr[22,0] = 0.9595582875050894
$bbArr[23,0] = 1.503966953576466
$bbArr[24,0] = 1.799837015295822
$bbArr[25,0] = 0.6518403676065248
$bbArr[26,0] = 0.6633823054011998
$bbArr[27,0] = 0.923139910942723
$bbArr[28,0] = 0.3869820931359413
$bbArr[29,0] = 2.093916965767463
$bbArr[30,0] = 1.580888475204972
$bbArr[31,0] = 0.3008043112709089
$bbArr[32,0] = 1.199598313222268
$bbArr[33,0] = 0.4774400648527148
$bbArr[34,0] = 1.507463254996111
$bbArr[35,0] = 0.2179116434425623
$bbArr[36,0] = 0.6120689161334525
$bbArr[37,0] = 2.489390679284554
$bbArr[38,0] = 0.5389418434166515
$bbArr[39,0] = 2.40118094791471
$bbArr[40,0] = 1.043009620608657
$bbArr[41,0] = 1.830928398766659
$bbArr[42,0] = -0.3674870133197601
$bbArr[43,0] = 1.62717758729876
$bbArr[44,0] = 1.117271732844245
$bbArr[45,0] = 0.6525147083449099
$bbArr[46,0] = 1.4
$bbArr[47,0] = -0.3
$bbArr[48,0] = -0.3
$bbArr[49,0] = 0.1
$bbArr[50,0] = -1.138880770453937
$bbArr[51,0] = -16.88491062648744
$bbArr[52,0] = 9.224715108933083
$bbArr[53,0] = 3.283355339827622
$bbArr[54,0] = 4.432584407022276
$bbArr[55,0] = 2.509693347214139
$bbArr[56,0] = -0.4381048169788073
$bbArr[57,0] = 4.106981763725997
$bbArr[58,0] = 0.3842995656585515
$bbArr[59,0] = 2.277966437795897
$bbArr[60,0] = 2.507553358214992
$bbArr[61,0] = -2.399190900254823
$bbArr[62,0] = -1.115644072253531
$bbArr[63,0] = 0.1689348086957096
$bbArr[64,0] = -1.814969742946232
$bbArr[65,0] = -1.586779238813989
$bbArr[66,0] = -0.2996177924633514
$bbArr[67,0] = 2.757652919539751
$bbArr[68,0] = -0.08601690538415596
$bbArr[69,0] = -0.6801011570971838
$bbArr[70,0] = 1.538981993999982
$bbArr[71,0] = 1.68501852020853
$bbArr[72,0] = 0.03331000006224372
$bbArr[73,0] = 0.03331000006224372
$bbArr[74,0] = 0.03331000006224372
$bbArr[75,0] = 0.03331000006224372
$bbArr[76,0] = 0.03331000006224372
$bbArr[77,0] = 0.03331000006224372
$bbArr[78,0] = 0.03331000006224372
$bbArr[79,0] = 0.03331000006224372
$bbArr[80,0] = 0.03331000006224372
$bbArr[81,0] = 0.03331000006224372
$bbArr[82,0] = 0.03331000006224372

$ws.Range("BB1:BB83").Value = $bbArr
